{"js": "// Goal (per the commit diff):\n//  - The paragraph that ends with the \"https://sb.digital/\" hyperlink no\n//    longer carries the hidden \"_GoBack\" bookmark.\n//  - Four additional blank (\"Garamond\") paragraphs are inserted between the\n//    \"Online Marketing and SEO Plan\" heading and the \"Google Analytics\"\n//    heading (there were 0 blank paragraphs between them before, now 4).\n//  - Both headings become centered, bold, single-underlined.\n//  - The \"_GoBack\" bookmark now sits on the last (blank) paragraph of the\n//    document, which is where Word leaves it after the user's last edit\n//    landed there.\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\n// Locate the anchor paragraphs by their text so this is resilient to the\n// exact paragraph index.\nlet sbDigitalIndex = -1;\nlet onlineMarketingIndex = -1;\nlet googleAnalyticsIndex = -1;\nfor (let i = 0; i < items.length; i++) {\n  const t = items[i].text;\n  if (t.indexOf(\"https://sb.digital/\") !== -1) {\n    sbDigitalIndex = i;\n  } else if (t === \"Online Marketing and SEO Plan\") {\n    onlineMarketingIndex = i;\n  } else if (t === \"Google Analytics\") {\n    googleAnalyticsIndex = i;\n  }\n}\n\nif (sbDigitalIndex === -1 || onlineMarketingIndex === -1 || googleAnalyticsIndex === -1) {\n  throw new Error(\"Could not locate the expected anchor paragraphs.\");\n}\n\n// 1) Move the \"_GoBack\" bookmark off the sb.digital hyperlink paragraph; it\n// will be re-added to the final paragraph of the document below.\ncontext.document.deleteBookmark(\"_GoBack\");\n\n// 2) Insert four blank paragraphs between the two headings *before* the\n// headings are reformatted, so the new paragraphs inherit the plain\n// (\"Garamond\" font, no bold/underline/center) style that is currently on\n// the \"Online Marketing and SEO Plan\" paragraph, matching the other blank\n// paragraphs already in this section.\nconst onlineMarketingPara = items[onlineMarketingIndex];\nlet insertAfter = onlineMarketingPara;\nfor (let i = 0; i < 4; i++) {\n  insertAfter = insertAfter.insertParagraph(\"\", Word.InsertLocation.after);\n}\n\n// 3) Now format the two headings: centered, bold, single underline.\nonlineMarketingPara.alignment = Word.Alignment.centered;\nonlineMarketingPara.font.bold = true;\nonlineMarketingPara.font.underline = Word.UnderlineType.single;\n\nconst googleAnalyticsPara = items[googleAnalyticsIndex];\ngoogleAnalyticsPara.alignment = Word.Alignment.centered;\ngoogleAnalyticsPara.font.bold = true;\ngoogleAnalyticsPara.font.underline = Word.UnderlineType.single;\n\nawait context.sync();\n\n// 4) Put the \"_GoBack\" bookmark back on the very last paragraph of the body.\nconst allParagraphs = context.document.body.paragraphs;\nallParagraphs.load(\"text\");\nawait context.sync();\n\nconst lastParagraph = allParagraphs.items[allParagraphs.items.length - 1];\nconst lastRange = lastParagraph.getRange();\nlastRange.insertBookmark(\"_GoBack\");\n\nawait context.sync();\n", "ps1": "# Goal (per the commit diff):\n#  - The paragraph ending with the \"https://sb.digital/\" hyperlink no longer\n#    carries the hidden \"_GoBack\" bookmark.\n#  - Four additional blank (\"Garamond\") paragraphs are inserted between the\n#    \"Online Marketing and SEO Plan\" heading and the \"Google Analytics\"\n#    heading.\n#  - Both headings become centered, bold, single-underlined.\n#  - The \"_GoBack\" bookmark is re-added on the last (blank) paragraph of the\n#    document.\n\n$d = $word.ActiveDocument\n\n# Locate the anchor paragraphs by their text content so this does not\n# depend on a hard-coded paragraph index.\n$sbDigitalIndex = -1\n$onlineMarketingIndex = -1\n$googleAnalyticsIndex = -1\n\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n  $p = $d.Paragraphs.Item($i)\n  $t = $p.Range.Text\n  if ($t -match \"https://sb\\.digital/\") {\n    $sbDigitalIndex = $i\n  } elseif ($t -eq \"Online Marketing and SEO Plan`r\") {\n    $onlineMarketingIndex = $i\n  } elseif ($t -eq \"Google Analytics`r\") {\n    $googleAnalyticsIndex = $i\n  }\n}\n\nif ($sbDigitalIndex -eq -1 -or $onlineMarketingIndex -eq -1 -or $googleAnalyticsIndex -eq -1) {\n  throw \"Could not locate the expected anchor paragraphs.\"\n}\n\n# 1) Move the \"_GoBack\" bookmark off the sb.digital hyperlink paragraph; it\n# is re-added to the final paragraph of the document below.\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n  $d.Bookmarks.Item(\"_GoBack\").Delete()\n}\n\n# 2) Insert four blank paragraphs right after the \"Online Marketing and SEO\n# Plan\" heading, before that heading gets reformatted, so the new\n# paragraphs inherit its current plain (\"Garamond\" font, no\n# bold/underline/center) style -- matching the other blank paragraphs\n# already in this section.\n$onlineMarketingPara = $d.Paragraphs.Item($onlineMarketingIndex)\nfor ($i = 0; $i -lt 4; $i++) {\n  $onlineMarketingPara.Range.InsertParagraphAfter()\n}\n\n# Re-resolve the heading paragraphs (indices shifted because of the insert).\n$googleAnalyticsIndex = $googleAnalyticsIndex + 4\n$onlineMarketingPara = $d.Paragraphs.Item($onlineMarketingIndex)\n$googleAnalyticsPara = $d.Paragraphs.Item($googleAnalyticsIndex)\n\n# 3) Format both headings: centered, bold, single underline.\n$onlineMarketingPara.Alignment = [Microsoft.Office.Interop.Word.WdParagraphAlignment]::wdAlignParagraphCenter\n$onlineMarketingPara.Range.Bold = 1\n$onlineMarketingPara.Range.Font.Underline = [Microsoft.Office.Interop.Word.WdUnderline]::wdUnderlineSingle\n\n$googleAnalyticsPara.Alignment = [Microsoft.Office.Interop.Word.WdParagraphAlignment]::wdAlignParagraphCenter\n$googleAnalyticsPara.Range.Bold = 1\n$googleAnalyticsPara.Range.Font.Underline = [Microsoft.Office.Interop.Word.WdUnderline]::wdUnderlineSingle\n\n# 4) Put the \"_GoBack\" bookmark back on the very last paragraph of the\n# document.\n$lastParagraph = $d.Paragraphs.Item($d.Paragraphs.Count)\n$d.Bookmarks.Add(\"_GoBack\", $lastParagraph.Range)\n"}
